$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet refs
# ---------------------------------------------------------------------------
$wsUsers   = $wb.Worksheets.Item("Users")
$wsWeekly  = $wb.Worksheets.Item("WeeklyEntryMatrix")
$wsProject = $wb.Worksheets.Item("Project_Title")
$wsActivity= $wb.Worksheets.Item("Activity_List")
$wsTimer   = $wb.Worksheets.Item("Update_Timer")
$wsHours   = $wb.Worksheets.Item("Update_Hours")

# ---------------------------------------------------------------------------
# Users sheet: just move the cursor, autofit column A
# ---------------------------------------------------------------------------
$wsUsers.Activate() | Out-Null
$wsUsers.Columns("A").ColumnWidth = 15.608072916666666
$wsUsers.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# WeeklyEntryMatrix sheet: rename the outsourced-contractor project string,
# bold the header, widen column A, move the cursor
# ---------------------------------------------------------------------------
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A2").Value = "Bartush-Cotton Creek Capital Management LLC-FVA-110095"
$wsWeekly.Range("A1").Font.Bold = $true
$wsWeekly.Columns("A").ColumnWidth = 50.053385416666664
$wsWeekly.Range("A11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Project_Title sheet: rename the engagement project string, bold the
# header row, widen columns A and B, move the cursor. This is the sheet
# that ends up active/selected.
# ---------------------------------------------------------------------------
$wsProject.Activate() | Out-Null
$wsProject.Range("A2").Value = "Project Wildcat-Ares Management LLC-FVA-108747"
$wsProject.Range("A1:B1").Font.Bold = $true
$wsProject.Columns("A").ColumnWidth = 42.721354166666664
$wsProject.Columns("B").ColumnWidth = 14.721354166666666
$wsProject.Range("G14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Activity_List sheet: cursor stays put (no content/format change), but it
# is no longer the active tab once we move on.
# ---------------------------------------------------------------------------
$wsActivity.Activate() | Out-Null
$wsActivity.Range("D11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Update_Timer sheet: bold the header, move the cursor
# ---------------------------------------------------------------------------
$wsTimer.Activate() | Out-Null
$wsTimer.Range("A1:B1").Font.Bold = $true
$wsTimer.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Update_Hours sheet: bold the header, move the cursor
# ---------------------------------------------------------------------------
$wsHours.Activate() | Out-Null
$wsHours.Range("A1:B1").Font.Bold = $true
$wsHours.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Finish back on Project_Title, which is the sheet left active/selected.
# ---------------------------------------------------------------------------
$wsProject.Activate() | Out-Null
$wsProject.Range("G14").Select() | Out-Null
